# Add a new "Spain" worksheet to the workbook, based on a copy of "Italy",
# matching the Gallery_MainProcessor_Repeaters_P_Panels.xlsx commit that
# added Test data for the Spain Zettler market.

$wb = $excel.ActiveWorkbook

# Clear the current tab's own selection/active status before we touch anything
# else, so it ends up looking like a normal (non-active) tab afterwards.
$italy = $wb.Worksheets.Item("Italy")

# 1) Duplicate "Italy" and place the copy right after it, then rename it.
$italy.Copy([System.Reflection.Missing]::Value, $italy)
$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"

# 2) Update the market-specific content on the new sheet.
#    Insert B4's new shared string ("NGC-3103/T2050/") before B2's
#    ("Spain Market") so the shared-string table gets the same ordering
#    the original author's Excel session produced.
$spain.Range("B4").Value = "NGC-3103/T2050/"
$spain.Range("B2").Value = "Spain Market"

# 3) Re-size columns B and D on the new sheet (widened versus "Italy").
$spain.Columns.Item(2).ColumnWidth = 21.7
$spain.Columns.Item(4).ColumnWidth = 34.7

# 4) Rows 3-5 grow to a taller, wrapped-text row height on the new sheet.
$spain.Rows.Item(3).RowHeight = 28.8
$spain.Rows.Item(4).RowHeight = 28.8
$spain.Rows.Item(5).RowHeight = 28.8

# 5) Restore a "whole sheet" style selection on "Italy" (it is no longer the
#    active tab) and leave the new "Spain" sheet with its own selection,
#    active and the visible/selected tab - matching the diff's tab switch.
$italy.Range("A1:D21").Select()
$spain.Range("C8:C9").Select()
$spain.Activate()
